# Consolidate youth-related department labels ("Ambassador", "Pathfinder",
# "Young Adult") into a single "Youth" department. "Congregation" is left
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B33").Value = "Youth"
$ws.Range("B40").Value = "Youth"
$ws.Range("B41").Value = "Youth"
$ws.Range("B43").Value = "Youth"
$ws.Range("B44").Value = "Youth"

# Restore the view state (scroll position / active cell selection) recorded
# in the saved workbook.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$ws.Range("M40").Select() | Out-Null
